$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values, regenerated after switching from Strike# to K,
# recalculating std/mean, and writing s_vals. Row -> new K value.
$kValues = @{
    2 = 1; 3 = 0; 4 = 2; 5 = 1; 6 = 1; 7 = 2; 8 = 2; 9 = 1; 10 = 2;
    11 = 2; 12 = 1; 13 = 1; 14 = 1; 15 = 2; 16 = 1; 17 = 1; 18 = 1; 19 = 1; 20 = 0;
    21 = 2; 22 = 1; 23 = 4; 24 = 2; 25 = 1; 26 = 1; 27 = 1; 28 = 0; 29 = 0;
    30 = 3; 31 = 2; 32 = 1; 33 = 2; 34 = 2; 35 = 2; 36 = 1; 37 = 3; 38 = 0; 39 = 0;
    41 = 1; 42 = 1; 43 = 0; 44 = 0; 45 = 1; 46 = 1; 47 = 3; 48 = 2; 49 = 1;
    50 = 0; 51 = 1; 52 = 0; 53 = 1; 54 = 1; 55 = 1; 56 = 0; 57 = 1; 58 = 1; 59 = 3;
    60 = 1; 61 = 1; 62 = 2; 63 = 0; 64 = 1; 65 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
